$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the existing header style (AC1) into the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# New header labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record data for every player row (2-42): Wins=82, Losses=80, Ties=0
$ws.Range("AD2:AD42").Value = 82
$ws.Range("AE2:AE42").Value = 80
$ws.Range("AF2:AF42").Value = 0
